# Rework the Industry/Tickers table for rows 7-13.
#
# The underlying shared-string pool changes: the "Internet Retail" industry
# is moved up next to "Diagnostics Research"/"CHEK" with a new ticker
# (TKAT), the "Internet Content Information"/LKCO and the old
# "Internet Retail"/"OCG, RMBL, RMBL" rows are dropped, every remaining
# industry/ticker pair shifts up by one row, the "Exchange Traded Fund"
# tickers shrink from "CEZ, PXJ" to just "CEZ", and a brand new
# "Lodging"/"STAY" row is appended before "Insurance - Property Casualty".
#
# Read all source values first (via Value2, which preserves the exact
# underlying text/escaping) before overwriting any cells, so that moved
# values are never clobbered before they're copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origRow7Industry  = $ws.Cells.Item(7, 1).Value2   # "Biotechnology "
$origRow7Tickers   = $ws.Cells.Item(7, 2).Value2   # "RUBY, SLGL, OCUP"
$origRow8Industry  = $ws.Cells.Item(8, 1).Value2   # "Internet Retail "
$origRow10Industry = $ws.Cells.Item(10, 1).Value2  # "Telecom Services "
$origRow10Tickers  = $ws.Cells.Item(10, 2).Value2  # "SJR"
$origRow11Industry = $ws.Cells.Item(11, 1).Value2  # "Agricultural Inputs "
$origRow11Tickers  = $ws.Cells.Item(11, 2).Value2  # "UAN"
$origRow12Industry = $ws.Cells.Item(12, 1).Value2  # "Oil Gas E&P "
$origRow12Tickers  = $ws.Cells.Item(12, 2).Value2  # "TPL"
$origRow13Industry = $ws.Cells.Item(13, 1).Value2  # "Exchange Traded Fund "

# Row 7: Internet Retail / TKAT (new ticker)
$ws.Cells.Item(7, 1).Value = $origRow8Industry
$ws.Cells.Item(7, 2).Value = "TKAT"

# Row 8: Biotechnology / RUBY, SLGL, OCUP (moved from old row 7)
$ws.Cells.Item(8, 1).Value = $origRow7Industry
$ws.Cells.Item(8, 2).Value = $origRow7Tickers

# Row 9: Telecom Services / SJR (moved from old row 10)
$ws.Cells.Item(9, 1).Value = $origRow10Industry
$ws.Cells.Item(9, 2).Value = $origRow10Tickers

# Row 10: Agricultural Inputs / UAN (moved from old row 11)
$ws.Cells.Item(10, 1).Value = $origRow11Industry
$ws.Cells.Item(10, 2).Value = $origRow11Tickers

# Row 11: Oil Gas E&P / TPL (moved from old row 12)
$ws.Cells.Item(11, 1).Value = $origRow12Industry
$ws.Cells.Item(11, 2).Value = $origRow12Tickers

# Row 12: Exchange Traded Fund / CEZ (industry moved from old row 13, tickers trimmed)
$ws.Cells.Item(12, 1).Value = $origRow13Industry
$ws.Cells.Item(12, 2).Value = "CEZ"

# Row 13: Lodging / STAY (brand new row)
$ws.Cells.Item(13, 1).Value = "Lodging "
$ws.Cells.Item(13, 2).Value = "STAY"
